# Progress update as of 04-Nov-2025: the "Period to expire" counters move by
# one day, and the "Last update" date advances from 03-Nov-2025 to
# 04-Nov-2025 on the Training Dashboard sheet (rows 3 and 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Row 3 (Endangered by Electricity A safety Training (SOPs)) ---
$ws.Range("H3").Value = -99

# Write the date as a literal text formula first, then paste the computed
# value back over itself (values only). This keeps I3 a plain text cell
# ("04-Nov-2025"), matching how the existing "LAST UPDATE" column is
# stored, instead of Excel's automatic text -> date-serial coercion that a
# plain Value/Value2 assignment of a date-shaped string would trigger.
$ws.Range("I3").Formula = '="04-Nov-2025"'
$ws.Range("I3").Copy()
$ws.Range("I3").PasteSpecial(-4163)

# --- Row 4 (IS0 55001 (Other Trainings)) ---
$ws.Range("H4").Value = 286

$ws.Range("I4").Formula = '="04-Nov-2025"'
$ws.Range("I4").Copy()
$ws.Range("I4").PasteSpecial(-4163)
